$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.446.66"
$ws.Range("E2").Value = "  -5.70%  "
$ws.Range("D3").Value = "'2.897.15"
$ws.Range("E3").Value = "  -3.52%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'550.63"
$ws.Range("E5").Value = "  -2.26%  "
$ws.Range("D6").Value = "'122.80"
$ws.Range("E6").Value = "  -4.68%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").Value = "'2.891.68"
$ws.Range("E8").Value = "  -3.67%  "
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("E10").Value = "  -8.27%  "
$ws.Range("D11").Value = "'4.71"
$ws.Range("E11").Value = "  -10.15%  "
$ws.Range("D12").Value = "'0.437"
$ws.Range("E12").Value = "  +1.55%  "
$ws.Range("E13").Value = "  -5.34%  "
$ws.Range("D14").Value = "'32.37"
$ws.Range("E14").Value = "  -1.61%  "
$ws.Range("E15").Value = "  +0.29%  "
$ws.Range("D16").Value = "'3.373.68"
$ws.Range("E16").Value = "  -3.50%  "
$ws.Range("D17").Value = "'2.895.34"
$ws.Range("E17").Value = "  -3.38%  "
$ws.Range("E18").Value = "  +5.19%  "
$ws.Range("D19").Value = "'57.430.85"
$ws.Range("E19").Value = "  -5.92%  "
$ws.Range("D20").Value = "'404.61"
$ws.Range("E20").Value = "  -7.65%  "
$ws.Range("D21").Value = "'12.89"
$ws.Range("E21").Value = "  -2.05%  "
$ws.Range("D22").Value = "'0.669"
$ws.Range("E22").Value = "  +0.74%  "
$ws.Range("D23").Value = "'6.82"
$ws.Range("E23").Value = "  -4.66%  "
$ws.Range("D24").Value = "'12.78"
$ws.Range("E24").Value = "  +1.67%  "
$ws.Range("D25").Value = "'76.90"
$ws.Range("E25").Value = "  -2.71%  "
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E28").Value = "  -1.83%  "
$ws.Range("E29").Value = "  +1.97%  "
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("D31").Value = "'6.02"
$ws.Range("E31").Value = "  -2.30%  "
$ws.Range("D32").Value = "'24.68"
$ws.Range("E32").Value = "  -3.31%  "
$ws.Range("D33").Value = "'0.0974"
$ws.Range("E33").Value = "  +3.66%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'5.42"
$ws.Range("E34").Value = "  -2.31%  "
$ws.Range("B35").Value = "Mantle"
$ws.Range("C35").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D35").Value = "'0.906"
$ws.Range("E35").Value = "  -4.47%  "
$ws.Range("E36").Value = "  -11.95%  "
$ws.Range("D37").Value = "'47.91"
$ws.Range("E37").Value = "  -4.58%  "
$ws.Range("D38").Value = "'8.36"
$ws.Range("E38").Value = "  +7.75%  "
$ws.Range("D39").Value = "'0.0₃0618"
$ws.Range("E40").Value = "  -1.52%  "
$ws.Range("E41").Value = "  -5.88%  "
$ws.Range("D42").Value = "'2.618.09"
$ws.Range("E42").Value = "  -1.85%  "
$ws.Range("D43").Value = "'360.04"
$ws.Range("E43").Value = "  -3.87%  "
$ws.Range("E44").Value = "  -2.50%  "
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").Value = "'119.02"
$ws.Range("E46").Value = "  +0.47%  "
$ws.Range("D47").Value = "'0.229"
$ws.Range("E47").Value = "  -3.03%  "
$ws.Range("D48").Value = "'0.107"
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("D49").Value = "'1.94"
$ws.Range("E49").Value = "  -1.58%  "
$ws.Range("D50").Value = "'22.87"
$ws.Range("E50").Value = "  -2.40%  "
$ws.Range("D51").Value = "'1.94"
$ws.Range("E51").Value = "  -4.63%  "
